$d = $word.ActiveDocument

# The first paragraph ("Chapter <del>10</del><ins>10</ins>: Glossary") has its
# "_GoBack" bookmark sitting right after the tracked-change insertion. Move it
# to the very start of the paragraph (right after the paragraph properties,
# before the "Chapter " run) as in the target revision. We rebuild the whole
# paragraph's run content via InsertXML (preserving every existing attribute)
# and then restore the Heading1 paragraph style via the object model, since
# setting the style that way is what correctly re-emits <w:pStyle>.

$savedTrack = $d.TrackRevisions
$d.TrackRevisions = $false

$r = $d.Paragraphs(1).Range
$r.InsertXML(@'
<w:p w14:paraId="07B36EBB" w14:textId="7D92B936" w:rsidR="000D2D15" w:rsidRDefault="00FC4A42" w:rsidP="00DE180B"><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve">Chapter </w:t></w:r><w:del w:id="1" w:author="Greg Landry" w:date="2017-03-09T10:32:00Z"><w:r w:rsidDel="009961C3"><w:delText>10</w:delText></w:r></w:del><w:ins w:id="2" w:author="Greg Landry" w:date="2017-03-09T10:35:00Z"><w:r w:rsidR="007C1EFC"><w:t>10</w:t></w:r></w:ins><w:r w:rsidR="008533BE"><w:t>: Glossary</w:t></w:r></w:p>
'@)

$d.Paragraphs(1).Style = "Heading1"

$d.TrackRevisions = $savedTrack
